$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for all existing data rows
# (rows 2 through 132) from 2023-09-21 (45190) to 2023-09-23 (45192).
for ($r = 2; $r -le 132; $r++) {
    $ws.Cells.Item($r, 3).Value = 45192
}

# Row 132 gains an explicit row height (15, custom height) in the new file.
$ws.Rows.Item(132).RowHeight = 15

# Add the new record as row 133.
$row = 133

$ws.Cells.Item($row, 1).Value = "A 45154-2023"

$ws.Cells.Item($row, 2).Value = 45191
$ws.Cells.Item($row, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item($row, 3).Value = 45192
$ws.Cells.Item($row, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item($row, 4).Value = "VÄRMLANDS LÄN"
$ws.Cells.Item($row, 5).Value = "STORFORS"

$ws.Cells.Item($row, 7).Value = 1
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = 0
$ws.Cells.Item($row, 15).Value = 0
$ws.Cells.Item($row, 16).Value = 0
$ws.Cells.Item($row, 17).Value = 0

$ws.Cells.Item($row, 18).Value = ""
$ws.Cells.Item($row, 18).WrapText = $true
